$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 -- shifts the existing rows 23:42 down to 24:43,
# matching the target diff (new record inserted, everything below pushed down
# by one row).
$ws.Rows.Item(23).Insert()

# Fill in the newly inserted row 23 with the new record's data. The columns
# A, B, C, E, F, G, I, R are constant for this market/category/quality block,
# same as every other row in the table.
$ws.Cells.Item(23, 1).Value  = 7
$ws.Cells.Item(23, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value  = "Ñuble"
$ws.Cells.Item(23, 4).Value  = 44539
$ws.Cells.Item(23, 5).Value  = 16
$ws.Cells.Item(23, 6).Value  = 100112021
$ws.Cells.Item(23, 7).Value  = "Ají"
$ws.Cells.Item(23, 8).Value  = "Americana (o)"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 15000
$ws.Cells.Item(23, 12).Value = 16000
$ws.Cells.Item(23, 13).Value = 15500
$ws.Cells.Item(23, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(23, 15).Value = "Región del Maule"
$ws.Cells.Item(23, 16).Value = 1033
$ws.Cells.Item(23, 17).Value = 15
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of the
# table (same format as D22 / D24, etc.).
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat
